# Updates cryptos list data (price + volume columns) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''28.165.56'
$ws.Range("E2").Value = '  -3.15%  '

# Row 3
$ws.Range("D3").Value = '''1.918.15'
$ws.Range("E3").Value = '  -3.84%  '

# Row 4
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -1.34%  '

# Row 5
$ws.Range("D5").Value = '''327.66'
$ws.Range("E5").Value = '  -0.81%  '

# Row 6
$ws.Range("D6").Value = '''1.002'
$ws.Range("E6").Value = '  -1.08%  '

# Row 7
$ws.Range("D7").Value = '''0.4679'
$ws.Range("E7").Value = '  -5.72%  '

# Row 8
$ws.Range("D8").Value = '''0.4016'
$ws.Range("E8").Value = '  -4.00%  '

# Row 9
$ws.Range("D9").Value = '''53.01'
$ws.Range("E9").Value = '  -3.59%  '

# Row 10
$ws.Range("D10").Value = '''0.08412'
$ws.Range("E10").Value = '  -5.09%  '

# Row 11
$ws.Range("D11").Value = '''1.046'
$ws.Range("E11").Value = '  -4.20%  '

# Row 12
$ws.Range("D12").Value = '''22.16'
$ws.Range("E12").Value = '  -3.02%  '

# Row 13
$ws.Range("D13").Value = '''1.922.63'
$ws.Range("E13").Value = '  -3.87%  '

# Row 14
$ws.Range("D14").Value = '''7.439'
$ws.Range("E14").Value = '  -6.68%  '

# Row 15
$ws.Range("D15").Value = '''6.069'
$ws.Range("E15").Value = '  -5.37%  '

# Row 16
$ws.Range("D16").Value = '''1.002'
$ws.Range("E16").Value = '  -1.39%  '

# Row 17
$ws.Range("D17").Value = '''89.65'
$ws.Range("E17").Value = '  -3.01%  '

# Row 18
$ws.Range("D18").Value = '''0.00001069'
$ws.Range("E18").Value = '  -3.21%  '

# Row 19
$ws.Range("D19").Value = '''0.06605'
$ws.Range("E19").Value = '  -2.14%  '

# Row 20
$ws.Range("D20").Value = '''17.93'
$ws.Range("E20").Value = '  -7.75%  '

# Row 22
$ws.Range("D22").Value = '''5.741'
$ws.Range("E22").Value = '  -3.72%  '

# Row 23
$ws.Range("D23").Value = '''28.168.54'
$ws.Range("E23").Value = '  -3.29%  '

# Row 24
$ws.Range("E24").Value = '  -6.22%  '

# Row 25
$ws.Range("D25").Value = '''2.298'
$ws.Range("E25").Value = '  +0.26%  '

# Row 26
$ws.Range("D26").Value = '''2.138.26'
$ws.Range("E26").Value = '  -4.57%  '

# Row 27
$ws.Range("D27").Value = '''153.28'
$ws.Range("E27").Value = '  -2.36%  '

# Row 28
$ws.Range("D28").Value = '''20.04'
$ws.Range("E28").Value = '  -3.55%  '

# Row 29
$ws.Range("D29").Value = '''5.776'
$ws.Range("E29").Value = '  -7.74%  '

# Row 30
$ws.Range("D30").Value = '''2.139'
$ws.Range("E30").Value = '  -4.66%  '

# Row 31
$ws.Range("D31").Value = '''123.54'
$ws.Range("E31").Value = '  -2.73%  '

# Row 32
$ws.Range("D32").Value = '''0.9777'

# Row 33
$ws.Range("D33").Value = '''0.09666'
$ws.Range("E33").Value = '  -1.98%  '

# Row 34
$ws.Range("D34").Value = '''1.440'
$ws.Range("E34").Value = '  -5.72%  '

# Row 35
$ws.Range("D35").Value = '''3.645'
$ws.Range("E35").Value = '  -2.31%  '

# Row 36
$ws.Range("D36").Value = '''5.548'
$ws.Range("E36").Value = '  -4.68%  '

# Row 37
$ws.Range("D37").Value = '''8.865'
$ws.Range("E37").Value = '  -2.17%  '

# Row 38
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '''1.266'
$ws.Range("E38").Value = '  -3.49%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.02299'
$ws.Range("E39").Value = '  -4.63%  '

# Row 40
$ws.Range("D40").Value = '''0.06181'
$ws.Range("E40").Value = '  -2.93%  '

# Row 41
$ws.Range("D41").Value = '''0.6172'
$ws.Range("E41").Value = '  -4.49%  '

# Row 42
$ws.Range("D42").Value = '''11.05'
$ws.Range("E42").Value = '  -4.24%  '

# Row 43
$ws.Range("E43").Value = '  -1.08%  '

# Row 44
$ws.Range("E44").Value = '  -3.46%  '

# Row 45
$ws.Range("D45").Value = '''1.315'
$ws.Range("E45").Value = '  -3.26%  '

# Row 46
$ws.Range("D46").Value = '''0.5856'
$ws.Range("E46").Value = '  -5.10%  '

# Row 47
$ws.Range("E47").Value = '  -4.61%  '

# Row 48
$ws.Range("D48").Value = '''2.026'
$ws.Range("E48").Value = '  -6.40%  '

# Row 49
$ws.Range("E49").Value = '  -1.56%  '

# Row 50
$ws.Range("D50").Value = '''0.06908'
$ws.Range("E50").Value = '  -0.61%  '

# Row 51
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '''0.00000000309'
$ws.Range("E51").Value = '  -11.24%  '
